# Updated code for serial run issue
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2 (User_TC001): clarify step 2 wording ---
$row2B = @"
To validate invite via email based on following conditions.
1) Launch application and login as global admin.
2) Click on Users Menu and customize columns
3) Click on Add User button.
4) Enter valid details in required field (Full Name, Email).
5) Select any one role (Carrier, Shipper Admin, Shipper User, Driver (Full) and Driver (Limited) .
6) Select any one value from Carrier Drop Down field.
7) Click on Invite button and check invitation email sent successfully.
8.) Clickon Resend Invitation and Check if invitation sent again.
"@
$ws.Range("B2").Value = $row2B

# --- Row 7 (Loads_TC001): drop periods after step letters, mark automated (YES) ---
$row7B = @"
Validate whehter Global admin is able to add new load using Shipper platform on following conditions.
a) Launch and login applcation as Global admin
b) Click on add new load button.
C) Set Carrier Name, Load Date, Shipper, Rate, Rate UOM and Commodity.
d) click on save and search for record in AG grid.
e) click on Edit and set Shipper contact and click on save.
f) Search for edited record in AG grid
g) Click on delete button.
h) Check whether record exist in AG grid after delete.
"@
$ws.Range("B7").Value = $row7B
$ws.Range("C7").Value = "YES"

# --- Row 8 (Loads_TC002): add trailing period, mark automated (YES) ---
$row8B = @"
Validate whehter carrier user can add Scoular loads for payment using full submit.
1) Enter valid user id and Password and click Login button.
2) Click on Add New Load button from Load menu.
3) Enter valid details in all required field and click Save button
4) Now loads are saved successfully.
5) Upload an Origin and Destination ticket image or PDF document for corresponding load.
6) Observe Ready to Submit Load icon in grid should change to green color.
7) Click on Submit Load button.
8) Select any option and click Submit button.

"@
$ws.Range("B8").Value = $row8B
$ws.Range("C8").Value = "YES"

# --- Row 9 (Loads_TC003): mark automated (YES) ---
$ws.Range("C9").Value = "YES"

# --- Row 11 (new Loads_TC005 test case) ---
$row11B = @"
Validate the Dispatch Schedule loads to respective Carrier.
1) Enter valid user id and Password and click Login button.
2) Click on Add New Load button from Load menu.
3) Enter valid details in required field and click Dispatch button.
4) Now Loads are Dispatched to corresponding Carrier.
"@

$ws.Range("A11").Value = "Loads_TC005"
$ws.Range("B11").Value = $row11B
$ws.Range("C11").Value = "YES"
$ws.Range("D11").Value = "Dispatch sent successfully and carrier has control over load."

$ws.Range("A11").VerticalAlignment = -4108
$ws.Range("C11").VerticalAlignment = -4108
$ws.Range("B11").WrapText = $true
$ws.Range("D11").WrapText = $true
$ws.Range("D11").VerticalAlignment = -4108

$ws.Rows.Item(11).RowHeight = 90

$ws.Range("C11").Select() | Out-Null
